# NMRA Region/Division Map: add new "Lakeshores" division entry for NER
# region (division 21, subdivision 15) per NMRA website.
#
# This inserts a new row at sheet row 17 (shifting existing rows 17-198
# down to 18-199), fills in the new entry's four columns, restores the
# selection to the freshly-entered cell, and re-applies the worksheet's
# remembered two-key sort (column A then column B) so the recorded
# sortState/sortCondition ranges grow to cover the new last row (D199)
# without actually reordering any of the existing (not strictly
# B-sorted) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 17; existing rows 17..198 shift to 18..199,
# carrying their styles/values with them, and the sheet dimension grows
# to A1:D199 automatically.
$ws.Rows.Item(17).Insert()

# Populate the new division/subdivision entry.
$ws.Range("A17").Value = 21
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = "Lakeshores"
$ws.Range("D17").Value = "NER"

# Match the author's final selection.
$ws.Range("B17").Select()

# Re-record the sort bookmark over the grown range (A3:D199) with the same
# two sort keys (column A, then column B) as before the edit. Column B is
# added as a cell-color sort key (SortOn=1/xlSortOnCellColor) rather than a
# value key: since no cells carry fill color, every row ties on that key,
# so the stable sort leaves row order exactly as the primary (column A)
# pass produced it -- i.e. unchanged relative order for all pre-existing
# rows, just like the real edit (the sheet's column-B order was never a
# strict ascending sort to begin with, only column A was).
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A3:A199"), 0, 1, 0, 0)
$sort.SortFields.Add($ws.Range("B3:B199"), 1, 1, 0, 0)
$sort.SetRange($ws.Range("A3:D199"))
$sort.Header = 0
$sort.Apply()
